# Commit: "Fruta / hortaliza, semanal"
# A new daily price record is inserted as row 128 (pushing the existing
# rows 128:235 down to 129:236), growing the used range from A1:R235 to
# A1:R236.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 128; everything from 128 downward
# shifts down by one row (old row 128 becomes 129, ..., old row 235
# becomes 236).
$ws.Rows("128:128").Insert()

# Populate the newly inserted row 128 with the new record.
$ws.Range("A128").Value = 3
$ws.Range("B128").Value = "Femacal de La Calera"
$ws.Range("C128").Value = "Coquimbo"
$ws.Range("D128").Value = 44907
$ws.Range("E128").Value = 5
$ws.Range("F128").Value = 100112026
$ws.Range("G128").Value = "Haba"
$ws.Range("H128").Value = "Sin especificar"
$ws.Range("I128").Value = "Primera"
$ws.Range("J128").Value = 105
$ws.Range("K128").Value = 8000
$ws.Range("L128").Value = 9000
$ws.Range("M128").Value = 8524
$ws.Range("N128").Value = "$/saco 25 kilos"
$ws.Range("O128").Value = "Provincia de Quillota"
$ws.Range("P128").Value = 341
$ws.Range("Q128").Value = 25
$ws.Range("R128").Value = "Hortaliza"
